$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.350.48'
$ws.Range('E2').Value = '  +2.50%  '
$ws.Range('D3').Value = '3.234.10'
$ws.Range('E3').Value = '  +4.78%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''575.91'
$ws.Range('E5').Value = '  +1.84%  '
$ws.Range('D6').Value = '''154.51'
$ws.Range('E6').Value = '  +8.45%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.223.32'
$ws.Range('E8').Value = '  +4.83%  '
$ws.Range('E9').Value = '  +3.67%  '
$ws.Range('D10').Value = '''7.08'
$ws.Range('E10').Value = '  +9.83%  '
$ws.Range('E11').Value = '  +4.79%  '
$ws.Range('E12').Value = '  +3.69%  '
$ws.Range('D13').Value = '''37.90'
$ws.Range('E13').Value = '  +5.85%  '
$ws.Range('E14').Value = '  +3.26%  '
$ws.Range('D15').Value = '3.752.38'
$ws.Range('E15').Value = '  +4.67%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = '''555.66'
$ws.Range('E16').Value = '  +12.30%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '66.369.85'
$ws.Range('E17').Value = '  +2.52%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '''0.115'
$ws.Range('E18').Value = '  +3.19%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.230.95'
$ws.Range('E19').Value = '  +4.28%  '
$ws.Range('E20').Value = '  +5.91%  '
$ws.Range('D21').Value = '''14.36'
$ws.Range('E21').Value = '  +4.20%  '
$ws.Range('D22').Value = '''0.739'
$ws.Range('E22').Value = '  +6.92%  '
$ws.Range('D23').Value = '''7.86'
$ws.Range('E23').Value = '  +9.30%  '
$ws.Range('D24').Value = '''13.60'
$ws.Range('E24').Value = '  +6.92%  '
$ws.Range('D25').Value = '''81.97'
$ws.Range('E25').Value = '  +3.86%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').Value = '''9.40'
$ws.Range('E27').Value = '  +18.03%  '
$ws.Range('E28').Value = '  +4.75%  '
$ws.Range('E29').Value = '  +8.21%  '
$ws.Range('D30').Value = '''27.75'
$ws.Range('E30').Value = '  +4.77%  '
$ws.Range('E31').Value = '  +2.42%  '
$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('E33').Value = '  +5.15%  '
$ws.Range('D34').Value = '''567.01'
$ws.Range('E34').Value = '  +10.11%  '
$ws.Range('D36').Value = '''6.43'
$ws.Range('E36').Value = '  +7.27%  '
$ws.Range('D37').Value = '''0.0462'
$ws.Range('E37').Value = '  +13.68%  '
$ws.Range('D38').Value = '''55.47'
$ws.Range('E38').Value = '  +3.75%  '
$ws.Range('D39').Value = '''0.0871'
$ws.Range('E39').Value = '  +8.76%  '
$ws.Range('D40').Value = '''3.04'
$ws.Range('E40').Value = '  +13.39%  '
$ws.Range('D41').Value = '''0.127'
$ws.Range('E41').Value = '  +4.43%  '
$ws.Range('D42').Value = '3.134.30'
$ws.Range('E42').Value = '  +6.57%  '
$ws.Range('D43').Value = '''8.62'
$ws.Range('E43').Value = '  +2.53%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').Value = '''0.274'
$ws.Range('E44').Value = '  +10.66%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = '''2.34'
$ws.Range('E45').Value = '  +7.63%  '
$ws.Range('D46').Value = '''27.07'
$ws.Range('E46').Value = '  +7.21%  '
$ws.Range('D47').Value = '0.0₃0562'
$ws.Range('E47').Value = '  +3.03%  '
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('E49').Value = '  +4.05%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '''122.54'
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').Value = '''2.25'
$ws.Range('E51').Value = '  +8.34%  '
